$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4027742743492126
$ws.Range("B1").Value = 0.4522325396537781
$ws.Range("C1").Value = 0.7532025575637817
$ws.Range("D1").Value = 1.560185194015503
$ws.Range("E1").Value = 2.850537061691284
